# Weekly fruit/vegetable update: insert a new daily price record as the
# most recent row (227) for Ajo (garlic) at Terminal Hortofrutícola Agro
# Chillán. All subsequent rows shift down by one (227->228, ..., 238->239).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 227; this pushes the existing
# rows 227-238 down to 228-239, carrying their values/styles with them.
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new weekly record.
$ws.Cells.Item(227, 1).Value2 = 7
$ws.Cells.Item(227, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(227, 3).Value = "Ñuble"
$ws.Cells.Item(227, 4).Value2 = 44753
$ws.Cells.Item(227, 5).Value2 = 16
$ws.Cells.Item(227, 6).Value2 = 100112003
$ws.Cells.Item(227, 7).Value = "Ajo"
$ws.Cells.Item(227, 8).Value = "Chino"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value2 = 80
$ws.Cells.Item(227, 11).Value2 = 20000
$ws.Cells.Item(227, 12).Value2 = 21000
$ws.Cells.Item(227, 13).Value2 = 20500
$ws.Cells.Item(227, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(227, 15).Value = "China"
$ws.Cells.Item(227, 16).Value2 = 2050
$ws.Cells.Item(227, 17).Value2 = 10
$ws.Cells.Item(227, 18).Value = "Hortaliza"
